# Refresh the cryptos list (Price + Volume(1h)) to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even if it looks numeric to Excel
# (e.g. "583.70" -> 583.7, "0.0000170" -> 1.7E-05). A leading apostrophe is
# Excel's standard quote-prefix text marker and is not stored as part of the
# cell value.
function Set-TextValue($addr, $text) {
    $ws.Range($addr).Value = "'" + $text
}

$ws.Range("D2").Value = '67.134.32'
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '2.487.68'
$ws.Range("E3").Value = '  -0.11%  '

$ws.Range("E4").Value = '  +0.08%  '

Set-TextValue "D5" '583.70'
$ws.Range("E5").Value = '  -0.45%  '

Set-TextValue "D6" '171.74'
$ws.Range("E6").Value = '  +2.15%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.97%  '

$ws.Range("D9").Value = '2.488.19'
$ws.Range("E9").Value = '  -0.08%  '

$ws.Range("E10").Value = '  +0.44%  '

$ws.Range("E11").Value = '  +0.12%  '

Set-TextValue "D12" '4.92'
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("E13").Value = '  -2.14%  '

Set-TextValue "D15" '25.37'
$ws.Range("E15").Value = '  -2.52%  '

$ws.Range("D16").Value = '67.084.19'
$ws.Range("E16").Value = '  +0.07%  '

Set-TextValue "D17" '0.0000170'
$ws.Range("E17").Value = '  -2.15%  '

$ws.Range("D18").Value = '2.484.48'
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("E19").Value = '  -6.29%  '

Set-TextValue "D20" '7.40'
$ws.Range("E20").Value = '  -5.09%  '

Set-TextValue "D21" '348.81'
$ws.Range("E21").Value = '  -3.38%  '

Set-TextValue "D22" '4.03'
$ws.Range("E22").Value = '  -2.07%  '

$ws.Range("E23").Value = '  -0.05%  '

Set-TextValue "D24" '68.57'
$ws.Range("E24").Value = '  -3.23%  '

$ws.Range("E25").Value = '  -4.92%  '

Set-TextValue "D26" '1.78'
$ws.Range("E26").Value = '  -3.46%  '

Set-TextValue "D27" '9.26'
$ws.Range("E27").Value = '  -2.11%  '

Set-TextValue "D28" '0.999'
$ws.Range("E28").Value = '  +0.48%  '

$ws.Range("D30").Value = '0.0₃0901'
$ws.Range("E30").Value = '  -3.88%  '

Set-TextValue "D31" '509.24'
$ws.Range("E31").Value = '  +0.97%  '

Set-TextValue "D32" '7.76'
$ws.Range("E32").Value = '  -4.19%  '

$ws.Range("E33").Value = '  -3.21%  '

$ws.Range("E34").Value = '  -4.15%  '

$ws.Range("E35").Value = '  +0.07%  '

Set-TextValue "D36" '159.77'
$ws.Range("E36").Value = '  +0.41%  '

$ws.Range("E37").Value = '  -7.92%  '

Set-TextValue "D38" '18.70'
$ws.Range("E38").Value = '  +0.72%  '

Set-TextValue "D39" '18.22'
$ws.Range("E39").Value = '  -5.17%  '

$ws.Range("E40").Value = '  -6.06%  '

$ws.Range("E41").Value = '  -2.66%  '

$ws.Range("E42").Value = '  -0.15%  '

$ws.Range("E43").Value = '  -2.31%  '

Set-TextValue "D44" '4.81'
$ws.Range("E44").Value = '  -3.22%  '

Set-TextValue "D45" '2.36'
$ws.Range("E45").Value = '  -4.66%  '

Set-TextValue "D46" '38.81'
$ws.Range("E46").Value = '  -1.32%  '

Set-TextValue "D47" '142.32'
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("E48").Value = '  -4.86%  '

Set-TextValue "D49" '3.44'
$ws.Range("E49").Value = '  -4.67%  '

$ws.Range("E50").Value = '  -6.45%  '

$ws.Range("E51").Value = '  -0.67%  '
